$wb = $excel.ActiveWorkbook

# --- Update "Hoja1" sheet: conversion text in A1 ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$wsHoja1.Range("A1").Value = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 12.45 = 50186.8 pesos`n✅ 50186.8 pesos = 12.42 = 974.54 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

# --- Update "tasas" sheet: rate values ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value = 80.3
$wsTasas.Range("O10").Value = 4030
$wsTasas.Range("N12").Value = 4040
$wsTasas.Range("O12").Value = 78.45
